# Updated for end of semester
# Add two new rows (Spring 2025) to the semester_reviews sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("semester_reviews")

# Row 8: Number=4, Instructor col=Course, Semester=Spring 2025
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Course"
$ws.Range("C8").Value = "Spring 2025"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 70

# Row 9: Number=4, Instructor col=Instructor, Semester=Spring 2025
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Instructor"
$ws.Range("C9").Value = "Spring 2025"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 70

# Update the selected cell to reflect where the user left off after entry.
$ws.Range("A11").Select()
